# Applies the Xhosa wording updates described in the commit:
#   - "uphononongo" ("the study") -> "uphando" ("the research") in a
#     handful of sentences, plus a couple of small rewordings.
# Uses Find (locate-only, wdReplaceNone) + direct Range.Text assignment
# so that literal straight/curly quote characters in the replacement
# text are preserved verbatim (Find/Replace's ReplaceWith argument
# otherwise gets passed through Word's AutoCorrect "smart quotes").

$d = $word.ActiveDocument

function Set-RangeText($find, $replace) {
    $rng = $d.Content
    $found = $rng.Find.Execute($find, $true, $false, $false, $false, $false, `
                                $true, 1, $false, $null, 0)
    if (-not $found) {
        throw "Text not found: $find"
    }
    $rng.Text = $replace
}

Set-RangeText `
    "Ukuba unayo nayiphi na imibuzo okanye iinkxalabo malunga namalungelo akho njengomthathi-nxaxheba kuphononongo, ungaqhagamshelana neqela lophononongo ku-" `
    "Ukuba unayo nayiphi na imibuzo okanye iinkxalabo malunga namalungelo akho njengomthathi-nxaxheba kuphando, ungaqhagamshelana neqela lophando ku-"

Set-RangeText `
    "Ukuba uneminye imibuzo okanye iinkxalabo malunga namalungelo akho, ungaqhagamshelana nenye yee komiti yokuziphatha edwelisiweyo: " `
    "Ukuba uneminye imibuzo okanye iinkxalabo malunga namalungelo akho, ungaqhagamshelana nenye yee komiti yemigaqo yokuziphatha edwelisiweyo: "

Set-RangeText `
    "Imvume eChaziweyo yokuThatha Inxaxheba kuPhononongo" `
    "Imvume eChaziweyo yokuThatha Inxaxheba kuPhando"

Set-RangeText `
    "Ndilufundile olu lwazi lungentla kwaye ndiyayazi into ekufuneka yenziwe ngumntwana wam." `
    "Ndilufundile olu lwazi lungentla kwaye ndiyayazi into elindelekileyo ngomntwana wam."

Set-RangeText `
    "Ungazigcina iinkcukacha zam zoqhagamshelwano zikhuselekile ukuze undixelele ngeziphumo zophononongo." `
    "Ungazigcina iinkcukacha zam zoqhagamshelwano zikhuselekile ukuze undixelele ngeziphumo zophando."

Set-RangeText `
    "Ukuba ulufundile kwaye waluqonda olu xwebhu lungasentla, uyavumelana nemiyalezo kwaye unike imvume yokuthatha inxaxheba kuphononongo, khetha u-" `
    "Ukuba ulufundile kwaye waluqonda olu xwebhu lungasentla, uyavumelana nemiyalezo kwaye unika imvume yokuthatha inxaxheba kuphando, khetha u-"

Set-RangeText `
    "Khetha `"Hayi`" kuWhatsApp ukuba awufuni ukuthatha inxaxheba." `
    "Khetha u `"Hayi`" kuWhatsApp ukuba awufuni ukuthatha inxaxheba."
